$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '245.78'
Set-TextValue $ws.Range('E2') '0.07%'
Set-TextValue $ws.Range('F2') '1-1-2023'
Set-TextValue $ws.Range('G2') '1'

# Row 3
Set-TextValue $ws.Range('D3') '25.93'
Set-TextValue $ws.Range('E3') '1.63%'
Set-TextValue $ws.Range('F3') '1-1-2023'
Set-TextValue $ws.Range('G3') '1'

# Row 4
Set-TextValue $ws.Range('D4') '5.164'
Set-TextValue $ws.Range('E4') '2.48%'
Set-TextValue $ws.Range('F4') '1-1-2023'
Set-TextValue $ws.Range('G4') '1'

# Row 5
Set-TextValue $ws.Range('D5') '0.05581'
Set-TextValue $ws.Range('E5') '-0.57%'
Set-TextValue $ws.Range('F5') '1-1-2023'
Set-TextValue $ws.Range('G5') '1'

# Row 6
Set-TextValue $ws.Range('D6') '6.473'
Set-TextValue $ws.Range('E6') '-1.40%'
Set-TextValue $ws.Range('F6') '1-1-2023'
Set-TextValue $ws.Range('G6') '1'

# Row 7
Set-TextValue $ws.Range('D7') '0.8140'
Set-TextValue $ws.Range('E7') '-0.63%'
Set-TextValue $ws.Range('F7') '1-1-2023'
Set-TextValue $ws.Range('G7') '1'

# Row 8
Set-TextValue $ws.Range('D8') '0.8427'
Set-TextValue $ws.Range('E8') '0.54%'
Set-TextValue $ws.Range('F8') '1-1-2023'
Set-TextValue $ws.Range('G8') '1'

# Row 9
Set-TextValue $ws.Range('D9') '0.06938'
Set-TextValue $ws.Range('E9') '-0.36%'
Set-TextValue $ws.Range('F9') '1-1-2023'
Set-TextValue $ws.Range('G9') '1'

# Row 10
Set-TextValue $ws.Range('D10') '0.02828'
Set-TextValue $ws.Range('E10') '0.23%'
Set-TextValue $ws.Range('F10') '1-1-2023'
Set-TextValue $ws.Range('G10') '1'

# Row 11
Set-TextValue $ws.Range('D11') '0.09387'
Set-TextValue $ws.Range('E11') '-0.14%'
Set-TextValue $ws.Range('F11') '1-1-2023'
Set-TextValue $ws.Range('G11') '1'

# Row 12
Set-TextValue $ws.Range('D12') '0.001507'
Set-TextValue $ws.Range('E12') '-0.34%'
Set-TextValue $ws.Range('F12') '1-1-2023'
Set-TextValue $ws.Range('G12') '1'

# Row 13
Set-TextValue $ws.Range('D13') '0.0006000'
Set-TextValue $ws.Range('E13') '0.83%'
Set-TextValue $ws.Range('F13') '1-1-2023'
Set-TextValue $ws.Range('G13') '1'

# Row 14
Set-TextValue $ws.Range('D14') '0.006199'
Set-TextValue $ws.Range('E14') '1.49%'
Set-TextValue $ws.Range('F14') '1-1-2023'
Set-TextValue $ws.Range('G14') '1'

# Row 15
Set-TextValue $ws.Range('D15') '3.607'
Set-TextValue $ws.Range('E15') '2.88%'
Set-TextValue $ws.Range('F15') '1-1-2023'
Set-TextValue $ws.Range('G15') '1'

# Row 16
Set-TextValue $ws.Range('D16') '3.019'
Set-TextValue $ws.Range('E16') '-0.05%'
Set-TextValue $ws.Range('F16') '1-1-2023'
Set-TextValue $ws.Range('G16') '1'

# Row 17
Set-TextValue $ws.Range('D17') '2.183'
Set-TextValue $ws.Range('E17') '4.37%'
Set-TextValue $ws.Range('F17') '1-1-2023'
Set-TextValue $ws.Range('G17') '1'

# Row 18
Set-TextValue $ws.Range('E18') '-2.11%'
Set-TextValue $ws.Range('F18') '1-1-2023'
Set-TextValue $ws.Range('G18') '1'

# Row 19
Set-TextValue $ws.Range('D19') '0.1331'
Set-TextValue $ws.Range('E19') '-0.70%'
Set-TextValue $ws.Range('F19') '1-1-2023'
Set-TextValue $ws.Range('G19') '1'

# Row 20
Set-TextValue $ws.Range('D20') '0.03124'
Set-TextValue $ws.Range('E20') '-3.32%'
Set-TextValue $ws.Range('F20') '1-1-2023'
Set-TextValue $ws.Range('G20') '1'

# Row 21
Set-TextValue $ws.Range('D21') '0.1272'
Set-TextValue $ws.Range('E21') '-4.91%'
Set-TextValue $ws.Range('F21') '1-1-2023'
Set-TextValue $ws.Range('G21') '1'

# Row 22
Set-TextValue $ws.Range('D22') '3.753'
Set-TextValue $ws.Range('E22') '-0.39%'
Set-TextValue $ws.Range('F22') '1-1-2023'
Set-TextValue $ws.Range('G22') '1'

# Row 23
Set-TextValue $ws.Range('D23') '0.04618'
Set-TextValue $ws.Range('E23') '-1.91%'
Set-TextValue $ws.Range('F23') '1-1-2023'
Set-TextValue $ws.Range('G23') '1'

# Row 24
Set-TextValue $ws.Range('E24') '2.46%'
Set-TextValue $ws.Range('F24') '1-1-2023'
Set-TextValue $ws.Range('G24') '1'

# Row 25
Set-TextValue $ws.Range('D25') '0.001249'
Set-TextValue $ws.Range('E25') '0.38%'
Set-TextValue $ws.Range('F25') '1-1-2023'
Set-TextValue $ws.Range('G25') '1'

# Row 26
Set-TextValue $ws.Range('D26') '0.004532'
Set-TextValue $ws.Range('E26') '5.89%'
Set-TextValue $ws.Range('F26') '1-1-2023'
Set-TextValue $ws.Range('G26') '1'

# Row 27
Set-TextValue $ws.Range('D27') '0.00009600'
Set-TextValue $ws.Range('E27') '-1.07%'
Set-TextValue $ws.Range('F27') '1-1-2023'
Set-TextValue $ws.Range('G27') '1'

# Row 28
Set-TextValue $ws.Range('D28') '0.0001661'
Set-TextValue $ws.Range('E28') '19.93%'
Set-TextValue $ws.Range('F28') '1-1-2023'
Set-TextValue $ws.Range('G28') '1'

# Row 29
Set-TextValue $ws.Range('F29') '1-1-2023'
Set-TextValue $ws.Range('G29') '1'

# Row 30
Set-TextValue $ws.Range('F30') '1-1-2023'
Set-TextValue $ws.Range('G30') '1'

# Row 31
Set-TextValue $ws.Range('F31') '1-1-2023'
Set-TextValue $ws.Range('G31') '1'

# Row 32
Set-TextValue $ws.Range('F32') '1-1-2023'
Set-TextValue $ws.Range('G32') '1'

# Row 33
Set-TextValue $ws.Range('F33') '1-1-2023'
Set-TextValue $ws.Range('G33') '1'

# Row 34
Set-TextValue $ws.Range('F34') '1-1-2023'
Set-TextValue $ws.Range('G34') '1'

# Row 35
Set-TextValue $ws.Range('F35') '1-1-2023'
Set-TextValue $ws.Range('G35') '1'

# Row 36
Set-TextValue $ws.Range('F36') '1-1-2023'
Set-TextValue $ws.Range('G36') '1'

# Row 37
Set-TextValue $ws.Range('F37') '1-1-2023'
Set-TextValue $ws.Range('G37') '1'

# Row 38
Set-TextValue $ws.Range('F38') '1-1-2023'
Set-TextValue $ws.Range('G38') '1'

# Row 39
Set-TextValue $ws.Range('F39') '1-1-2023'
Set-TextValue $ws.Range('G39') '1'

# Row 40
Set-TextValue $ws.Range('D40') '0.03647'
Set-TextValue $ws.Range('E40') '-0.45%'
Set-TextValue $ws.Range('F40') '1-1-2023'
Set-TextValue $ws.Range('G40') '1'

# Row 41
Set-TextValue $ws.Range('D41') '0.1368'
Set-TextValue $ws.Range('E41') '30.53%'
Set-TextValue $ws.Range('F41') '1-1-2023'
Set-TextValue $ws.Range('G41') '1'

# Row 42
Set-TextValue $ws.Range('B42') 'KickToken'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range('D42') '0.006150'
Set-TextValue $ws.Range('E42') '-1.48%'
Set-TextValue $ws.Range('F42') '1-1-2023'
Set-TextValue $ws.Range('G42') '1'

# Row 43
Set-TextValue $ws.Range('B43') 'CEJI'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws.Range('D43') '0.002537'
Set-TextValue $ws.Range('E43') '-2.46%'
Set-TextValue $ws.Range('F43') '1-1-2023'
Set-TextValue $ws.Range('G43') '1'

# Row 44
Set-TextValue $ws.Range('D44') '0.008015'
Set-TextValue $ws.Range('E44') '-5.57%'
Set-TextValue $ws.Range('F44') '1-1-2023'
Set-TextValue $ws.Range('G44') '1'

# Row 45
Set-TextValue $ws.Range('D45') '0.00005340'
Set-TextValue $ws.Range('E45') '0.89%'
Set-TextValue $ws.Range('F45') '1-1-2023'
Set-TextValue $ws.Range('G45') '1'

# Row 46
Set-TextValue $ws.Range('E46') '-0.03%'
Set-TextValue $ws.Range('F46') '1-1-2023'
Set-TextValue $ws.Range('G46') '1'

# Row 47
Set-TextValue $ws.Range('E47') '-19.45%'
Set-TextValue $ws.Range('F47') '1-1-2023'
Set-TextValue $ws.Range('G47') '1'

# Row 48
Set-TextValue $ws.Range('D48') '0.002408'
Set-TextValue $ws.Range('E48') '19.39%'
Set-TextValue $ws.Range('F48') '1-1-2023'
Set-TextValue $ws.Range('G48') '1'

# Row 49
Set-TextValue $ws.Range('E49') '-0.03%'
Set-TextValue $ws.Range('F49') '1-1-2023'
Set-TextValue $ws.Range('G49') '1'

# Row 50
Set-TextValue $ws.Range('E50') '-0.03%'
Set-TextValue $ws.Range('F50') '1-1-2023'
Set-TextValue $ws.Range('G50') '1'

# Row 51
Set-TextValue $ws.Range('F51') '1-1-2023'
Set-TextValue $ws.Range('G51') '1'
